$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update column C ("Förändrad") for rows 2..440 from 45182 -> 45184
for ($r = 2; $r -le 440; $r++) {
    $ws.Cells.Item($r, 3).Value = 45184
}

# 2) Row 440 gains an explicit custom row height (ht="15" customHeight="1")
$ws.Rows.Item(440).RowHeight = 15

# 3) Append the new row 441 with the new logging notice
$ws.Cells.Item(441, 1).Value = "A 42642-2023"
$ws.Cells.Item(441, 2).Value = 45181
$ws.Cells.Item(441, 3).Value = 45184
$ws.Cells.Item(441, 4).Value = "UPPSALA LÄN"
$ws.Cells.Item(441, 5).Value = "TIERP"
$ws.Cells.Item(441, 7).Value = 3
$ws.Cells.Item(441, 8).Value = 0
$ws.Cells.Item(441, 9).Value = 0
$ws.Cells.Item(441, 10).Value = 0
$ws.Cells.Item(441, 11).Value = 0
$ws.Cells.Item(441, 12).Value = 0
$ws.Cells.Item(441, 13).Value = 0
$ws.Cells.Item(441, 14).Value = 0
$ws.Cells.Item(441, 15).Value = 0
$ws.Cells.Item(441, 16).Value = 0
$ws.Cells.Item(441, 17).Value = 0

# Match the date-format style used by the other B/C columns
$ws.Cells.Item(441, 2).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(441, 3).NumberFormat = "YYYY-MM-DD"
# Match the wrap-text style used by column R (Artnamn) on every row
$ws.Cells.Item(441, 18).WrapText = $true

$ws.Rows.Item(441).RowHeight = 15
